$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handback
#
# The handback run produced a target file (a.md), a handback package file
# (the de-de/zh-cn xliff) and a handback timestamp for each language sheet,
# and flipped the Status column from "Ready for handoff" to
# "Handed back: in sync with en-US".
# ---------------------------------------------------------------------------

$statusText = "Handed back: in sync with en-US"

$zhHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$deHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

$zhHandbackTime = "2016-08-15 20:31:38"
$deHandbackTime = "2016-08-15 20:31:45"

$targetFileName = "a.md"
$targetFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e4629a1e6c2df686dfd9b364929b7c7602f3e3b3/e2e/a.md"

function Update-LanguageSheet {
    param([string]$SheetName, [string]$HandbackFile, [string]$HandbackTime)

    $ws = $wb.Worksheets.Item($SheetName)

    # Status column (C) for both data rows.
    $ws.Range("C2").Value = $statusText
    $ws.Range("C3").Value = $statusText

    # Latest Target File (I): link to the translated markdown file.
    $ws.Hyperlinks.Add($ws.Range("I2"), $targetFileUrl, "", "", $targetFileName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("I3"), $targetFileUrl, "", "", $targetFileName) | Out-Null

    # Latest Handback File (J): the generated handback xliff package.
    $ws.Range("J2").Value = $HandbackFile
    $ws.Range("J3").Value = $HandbackFile

    # Latest Handback DateTime (K): when the handback report was generated.
    $ws.Range("K2").Value = $HandbackTime
    $ws.Range("K3").Value = $HandbackTime

    # Column widths grew to fit the new (longer) status/file-name content.
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(10).ColumnWidth = 40
}

Update-LanguageSheet "zh-cn" $zhHandbackFile $zhHandbackTime
Update-LanguageSheet "de-de" $deHandbackFile $deHandbackTime

# Overview sheet: the zh-cn/de-de status columns mirror the same shared
# status string, and widen to fit it as well.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527
